$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the data values in row 2 (new EM algorithm results) ---
$ws.Range("B2").Value = 762957.95475691999
$ws.Range("C2").Value = 229725.89961944899
$ws.Range("D2").Value = 59192.511157275803
$ws.Range("E2").Value = 18436.239720580299
$ws.Range("F2").Value = 6450.4833894153599

# --- Change the number format of the data cells from scientific to 0.000 ---
$ws.Range("B2:F4").NumberFormat = "0.000"

# --- Narrow columns B:C slightly (content now fits tighter) ---
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 9.6

# --- Update the selected cell shown when the sheet is active ---
$ws.Range("H14").Select()

# --- Chart updates ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)

# Show data labels to the right of each point
$series.HasDataLabels = 1
$dLbls = $series.DataLabels()
$dLbls.ShowValue = 1
$dLbls.Position = -4152

# Value axis: make visible, with plain integer tick labels next to axis
$valAx = $chart.Axes(2)
$valAx.Visible = 1
$valAx.TickLabelPosition = 4
$valAx.MajorTickMark = 4
$valAx.NumberFormat = "0"
$valAx.NumberFormatLinked = 0

# Category axis: cross value axis at 0
$catAx = $chart.Axes(1)
$catAx.CrossesAt = 0

Write-Host "edit complete"
